$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.926.44'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '3.301.08'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '185.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '576.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.601'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.129'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.17%  '
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = '3.871.82'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = '67.236.81'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '3.285.11'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '442.28'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +9.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('E21').Value = '  +2.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').Value = '3.432.77'
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.09'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.06%  '
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('E30').Value = '  +1.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.82'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.34'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.14%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.79'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('E36').Value = '  +5.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.60'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.49'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.47'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.784'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('D42').Value = '2.732.60'
$ws.Range('E42').Value = '  +2.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.44%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.19'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.01%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0671'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '329.12'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.62%  '
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.993'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.24'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.83%  '
